{"js": "// Move the leading \"2022\" to the front of the sentence:\n//   \"Fechas de la campa\u00f1a para constelaci\u00f3n de ori\u00f3n 2022: ...\"\n// becomes\n//   \"2022 Fechas de la campa\u00f1a para constelaci\u00f3n de ori\u00f3n: ...\"\n// Applies to every occurrence of this sentence in the document body.\n\nconst oldText =\n  \"Fechas de la campa\u00f1a para constelaci\u00f3n de ori\u00f3n 2022: 16-25 de enero, 14-23 de febrero, 14-24 de marzo\";\nconst newText =\n  \"2022 Fechas de la campa\u00f1a para constelaci\u00f3n de ori\u00f3n: 16-25 de enero, 14-23 de febrero, 14-24 de marzo\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Move the leading \"2022\" to the front of the sentence:\n#   \"Fechas de la campa\u00f1a para constelaci\u00f3n de ori\u00f3n 2022: ...\"\n# becomes\n#   \"2022 Fechas de la campa\u00f1a para constelaci\u00f3n de ori\u00f3n: ...\"\n# Applies to every occurrence of this sentence in the document.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"Fechas de la campa\u00f1a para constelaci\u00f3n de ori\u00f3n 2022: 16-25 de enero, 14-23 de febrero, 14-24 de marzo\"\n$find.Replacement.Text = \"2022 Fechas de la campa\u00f1a para constelaci\u00f3n de ori\u00f3n: 16-25 de enero, 14-23 de febrero, 14-24 de marzo\"\n\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue - keep searching the whole document\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2 - replace every match in the document, not just the first\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, `\n    $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n"}
